$d = $word.ActiveDocument

# Locate the three consecutive paragraphs that must be removed:
#   1. the blank paragraph right after "LOQ4209: Engenharia da Qualidade I (Requisito fraco)"
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
# Walk the Paragraphs collection and find the "Ver no Jupiter..." paragraph by its text.

$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $idx = $i
        break
    }
}

if ($idx -gt 0) {
    $prev = $d.Paragraphs.Item($idx - 1)   # the blank paragraph before it
    $next = $d.Paragraphs.Item($idx + 1)   # the "© 2020 ..." paragraph

    $deleteRange = $d.Range($prev.Range.Start, $next.Range.End)
    $deleteRange.Delete()
}
